$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K" = strikeouts) values for rows 2-13 and 15-16.
# Row 14's G value (0) is unchanged.
$gValues = @{
    2  = 1
    3  = 3
    4  = 0
    5  = 1
    6  = 1
    7  = 2
    8  = 1
    9  = 0
    10 = 1
    11 = 0
    12 = 2
    13 = 1
    15 = 0
    16 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
